$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: "SEO" -> "Accessibilité"; "sitecheckerpro" -> "lighthouse"
$ws.Range("A12").Value = "Accessibilité"
$ws.Range("F12").Value = "lighthouse"

# Row 13: add "black hat SEO" in C13
$ws.Range("C13").Value = "black hat SEO"

# Row 17: turn into a new Accessibilité row about "label sur social"
$ws.Range("A17").Value = "Accessibilité"
$ws.Range("B17").Value = "label sur social"
$ws.Range("F17").Value = "lighthouse"
$ws.Range("G17").Value = "x"

# Update active selection to E20
$ws.Range("E20").Select()
